$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Set all cell VALUES first, in row-major (document) order. The order in
#    which string values are first assigned controls the order they end up
#    in xl/sharedStrings.xml, so this must mirror the final layout:
#    row1 (C1), row2 (B2,C2), row3 (A3,B3,C3), row4..row9 (B only).
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "Contributor"

$ws.Range("B2").Value = "Added a changelog-file, `naltered the .gitignore according to `nthe Kicad standard, started adding `nannotations to the datasheets for navigation."
$ws.Range("C2").Value = "Doomn00b"

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A3").Value = 44784
$ws.Range("B3").Value = "Changed the project to use generic`nlinks instead of local links, merged footprint`nlibraries into one, merged symbol libraries into one,`nstarted rewriting V9990 datasheet with Open Sans font,`nfor future open source datasheet"
$ws.Range("C3").Value = "Doomn00b"

$ws.Range("B4").Value = "Changed V9990 symbol and footprint"
$ws.Range("B5").Value = "Replaced C3 THT capacitor with SMD."

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A6").Value = 44785
$ws.Range("B6").Value = "Replaced RGB-connector with real VGA conn from Amphenol, changed decoupling, put power-parts on their own sheet,"

$ws.Range("B7").Value = "Replaced big 14MHz oscillator with SMD one,"
$ws.Range("B8").Value = "Changed the glue logic to use Advanced 74xx versions, with far higher performance."
$ws.Range("B9").Value = "Removed Composite video and replaced with RGB Scart"

# ---------------------------------------------------------------------------
# 2. Apply formatting. The order in which *new* cell-format combinations are
#    first used controls the order they are appended to xl/styles.xml, so
#    create the word-wrap style (B2) before the "Calculation" style (C1).
# ---------------------------------------------------------------------------
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 75

$ws.Range("B3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 135

$wb.Styles.Add("Calculation") | Out-Null
$ws.Range("C1").Style = "Calculation"

# ---------------------------------------------------------------------------
# 3. Column sizing to fit the new content
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 30.8
$ws.Columns.Item(3).ColumnWidth = 10.5

# ---------------------------------------------------------------------------
# 4. Selection / active cell, matching the saved workbook view
# ---------------------------------------------------------------------------
$ws.Range("B10").Select() | Out-Null
